$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'308.45"
$ws.Range('E2').Value = "'0.07%"
$ws.Range('D3').Value = "'41.14"
$ws.Range('E3').Value = "'0.83%"
$ws.Range('D4').Value = "'5.172"
$ws.Range('E4').Value = "'0.79%"
$ws.Range('D5').Value = "'0.07683"
$ws.Range('E5').Value = "'0.78%"
$ws.Range('D6').Value = "'1.644"
$ws.Range('E6').Value = "'2.33%"
$ws.Range('E7').Value = "'1.17%"
$ws.Range('D8').Value = "'2.425"
$ws.Range('E8').Value = "'-1.78%"
$ws.Range('D9').Value = "'0.1233"
$ws.Range('E9').Value = "'10.00%"
$ws.Range('D10').Value = "'0.1823"
$ws.Range('E10').Value = "'2.47%"
$ws.Range('D11').Value = "'0.09240"
$ws.Range('E11').Value = "'0.88%"
$ws.Range('D12').Value = "'0.04230"
$ws.Range('E12').Value = "'-1.11%"
$ws.Range('D13').Value = "'0.1052"
$ws.Range('E13').Value = "'0.04%"
$ws.Range('D14').Value = "'0.001259"
$ws.Range('E14').Value = "'0.71%"
$ws.Range('D15').Value = "'0.005757"
$ws.Range('E15').Value = "'-0.09%"
$ws.Range('E16').Value = "'1,903.77%"
$ws.Range('D17').Value = "'3.345"
$ws.Range('E17').Value = "'-0.13%"
$ws.Range('D18').Value = "'4.317"
$ws.Range('E18').Value = "'1.60%"
$ws.Range('E19').Value = "'1.31%"
$ws.Range('D20').Value = "'7.385"
$ws.Range('E20').Value = "'11.09%"
$ws.Range('D21').Value = "'0.1382"
$ws.Range('E21').Value = "'1.23%"
$ws.Range('D22').Value = "'0.2712"
$ws.Range('E22').Value = "'-2.81%"
$ws.Range('D23').Value = "'0.04028"
$ws.Range('E23').Value = "'-1.26%"
$ws.Range('D24').Value = "'0.001269"
$ws.Range('E24').Value = "'2.64%"
$ws.Range('D25').Value = "'0.004226"
$ws.Range('E25').Value = "'2.86%"
$ws.Range('D26').Value = "'0.0001302"
$ws.Range('E26').Value = "'0.07%"
$ws.Range('D38').Value = "'0.02520"
$ws.Range('E38').Value = "'5.39%"
$ws.Range('D39').Value = "'0.05313"
$ws.Range('E39').Value = "'2.48%"
$ws.Range('D40').Value = "'0.007852"
$ws.Range('E40').Value = "'0.91%"
$ws.Range('E41').Value = "'1.06%"
$ws.Range('D42').Value = "'0.006671"
$ws.Range('E42').Value = "'-5.42%"
$ws.Range('D43').Value = "'0.001863"
$ws.Range('E43').Value = "'-4.55%"
$ws.Range('D44').Value = "'0.007994"
$ws.Range('E44').Value = "'0.60%"
$ws.Range('D45').Value = "'0.3070"
$ws.Range('E45').Value = "'-0.42%"
$ws.Range('D46').Value = "'0.00006723"
$ws.Range('E46').Value = "'-4.03%"
$ws.Range('D47').Value = "'0.00000000751"
$ws.Range('E47').Value = "'0.06%"
$ws.Range('D48').Value = "'0.2958"
$ws.Range('E48').Value = "'837.59%"
$ws.Range('D49').Value = "'0.003104"
$ws.Range('E49').Value = "'-26.11%"
$ws.Range('D50').Value = "'0.00002103"
$ws.Range('E50').Value = "'0.06%"
$ws.Range('D51').Value = "'0.0002003"
$ws.Range('E51').Value = "'0.06%"
